$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1836734693877551
$ws.Range("C2").Value = 0.5836734693877551
$ws.Range("J2").Value = 0.02448979591836735
$ws.Range("P2").Value = 0.1224489795918367
$ws.Range("S2").Value = 0.08571428571428572
$ws.Range("B3").Value = 0.01333333333333333
$ws.Range("C3").Value = 0.03333333333333333
$ws.Range("J3").Value = 0.02666666666666667
$ws.Range("P3").Value = 0.7066666666666667
$ws.Range("S3").Value = 0.22
$ws.Range("B6").Value = 0.04526748971193416
$ws.Range("D6").Value = 0.00823045267489712
$ws.Range("F6").Value = 0.07407407407407407
$ws.Range("J6").Value = 0.2139917695473251
$ws.Range("O6").Value = 0.01646090534979424
$ws.Range("Q6").Value = 0.139917695473251
$ws.Range("R6").Value = 0.06995884773662552
$ws.Range("S6").Value = 0.4320987654320987
$ws.Range("B7").Value = 0.06145251396648044
$ws.Range("D7").Value = 0.0335195530726257
$ws.Range("F7").Value = 0.0446927374301676
$ws.Range("J7").Value = 0.1564245810055866
$ws.Range("O7").Value = 0.0111731843575419
$ws.Range("Q7").Value = 0.1899441340782123
$ws.Range("R7").Value = 0.06145251396648044
$ws.Range("S7").Value = 0.441340782122905
$ws.Range("B8").Value = 0.09503239740820735
$ws.Range("D8").Value = 0.01511879049676026
$ws.Range("F8").Value = 0.0734341252699784
$ws.Range("J8").Value = 0.1058315334773218
$ws.Range("O8").Value = 0.02159827213822894
$ws.Range("Q8").Value = 0.1663066954643629
$ws.Range("R8").Value = 0.123110151187905
$ws.Range("S8").Value = 0.3995680345572354
$ws.Range("B9").Value = 0.0639269406392694
$ws.Range("D9").Value = 0.0273972602739726
$ws.Range("F9").Value = 0.0639269406392694
$ws.Range("J9").Value = 0.0821917808219178
$ws.Range("O9").Value = 0.0045662100456621
$ws.Range("Q9").Value = 0.1917808219178082
$ws.Range("R9").Value = 0.1095890410958904
$ws.Range("S9").Value = 0.45662100456621
$ws.Range("B10").Value = 0.09480626545754328
$ws.Range("D10").Value = 0.02555647155812036
$ws.Range("E10").Value = 0.001648804616652927
$ws.Range("F10").Value = 0.07831821929101401
$ws.Range("J10").Value = 0.112118713932399
$ws.Range("O10").Value = 0.01483924154987634
$ws.Range("Q10").Value = 0.1978565539983512
$ws.Range("R10").Value = 0.09563066776586975
$ws.Range("S10").Value = 0.3792250618301731
$ws.Range("G11").Value = 0.157088122605364
$ws.Range("J11").Value = 0.06130268199233716
$ws.Range("K11").Value = 0.210727969348659
$ws.Range("L11").Value = 0.5670498084291188
$ws.Range("S11").Value = 0.003831417624521073
$ws.Range("G12").Value = 0.7375
$ws.Range("J12").Value = 0.1875
$ws.Range("K12").Value = 0.00625
$ws.Range("L12").Value = 0.05625
$ws.Range("S12").Value = 0.0125
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.4047619047619048
$ws.Range("S13").Value = 0.02380952380952381
$ws.Range("F15").Value = 0.0154639175257732
$ws.Range("H15").Value = 0.1391752577319588
$ws.Range("I15").Value = 0.1082474226804124
$ws.Range("J15").Value = 0.3608247422680412
$ws.Range("K15").Value = 0.07731958762886598
$ws.Range("O15").Value = 0.06701030927835051
$ws.Range("S15").Value = 0.2319587628865979
$ws.Range("F16").Value = 0.01807228915662651
$ws.Range("H16").Value = 0.1506024096385542
$ws.Range("I16").Value = 0.07228915662650602
$ws.Range("J16").Value = 0.4819277108433735
$ws.Range("K16").Value = 0.108433734939759
$ws.Range("M16").Value = 0.02409638554216868
$ws.Range("N16").Value = 0.006024096385542169
$ws.Range("O16").Value = 0.04819277108433735
$ws.Range("S16").Value = 0.09036144578313253
$ws.Range("F17").Value = 0.0380952380952381
$ws.Range("H17").Value = 0.1880952380952381
$ws.Range("I17").Value = 0.09047619047619047
$ws.Range("J17").Value = 0.4214285714285714
$ws.Range("K17").Value = 0.09047619047619047
$ws.Range("M17").Value = 0.0119047619047619
$ws.Range("O17").Value = 0.05
$ws.Range("S17").Value = 0.1095238095238095
$ws.Range("F18").Value = 0.02232142857142857
$ws.Range("H18").Value = 0.1741071428571428
$ws.Range("I18").Value = 0.09821428571428571
$ws.Range("J18").Value = 0.4598214285714285
$ws.Range("K18").Value = 0.04910714285714286
$ws.Range("M18").Value = 0.02232142857142857
$ws.Range("O18").Value = 0.04910714285714286
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01579778830963665
$ws.Range("H19").Value = 0.2369668246445498
$ws.Range("I19").Value = 0.0995260663507109
$ws.Range("J19").Value = 0.3467614533965245
$ws.Range("K19").Value = 0.09873617693522907
$ws.Range("M19").Value = 0.02290679304897314
$ws.Range("O19").Value = 0.06161137440758294
$ws.Range("S19").Value = 0.1176935229067931
